$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 930-1016: shift data block down by 3 rows, new latest-week data at top ---
$ws.Cells.Item(930, 4).Value = 44769
$ws.Cells.Item(930, 9).Value = "Primera"
$ws.Cells.Item(930, 10).Value = 15000
$ws.Cells.Item(930, 11).Value = 150
$ws.Cells.Item(930, 12).Value = 150
$ws.Cells.Item(930, 13).Value = 150
$ws.Cells.Item(930, 16).Value = 150
$ws.Cells.Item(931, 4).Value = 44769
$ws.Cells.Item(931, 9).Value = "Segunda"
$ws.Cells.Item(931, 10).Value = 19000
$ws.Cells.Item(931, 11).Value = 120
$ws.Cells.Item(931, 12).Value = 120
$ws.Cells.Item(931, 13).Value = 120
$ws.Cells.Item(931, 16).Value = 120
$ws.Cells.Item(932, 4).Value = 44769
$ws.Cells.Item(932, 9).Value = "Tercera"
$ws.Cells.Item(932, 10).Value = 8000
$ws.Cells.Item(932, 11).Value = 80
$ws.Cells.Item(932, 12).Value = 80
$ws.Cells.Item(932, 13).Value = 80
$ws.Cells.Item(932, 16).Value = 80
$ws.Cells.Item(933, 4).Value = 44596
$ws.Cells.Item(933, 9).Value = "Primera"
$ws.Cells.Item(933, 10).Value = 59000
$ws.Cells.Item(933, 11).Value = 75
$ws.Cells.Item(933, 12).Value = 80
$ws.Cells.Item(933, 13).Value = 77
$ws.Cells.Item(933, 16).Value = 77
$ws.Cells.Item(934, 4).Value = 44596
$ws.Cells.Item(934, 9).Value = "Segunda"
$ws.Cells.Item(934, 10).Value = 49000
$ws.Cells.Item(934, 11).Value = 60
$ws.Cells.Item(934, 12).Value = 65
$ws.Cells.Item(934, 13).Value = 62
$ws.Cells.Item(934, 16).Value = 62
$ws.Cells.Item(935, 4).Value = 44596
$ws.Cells.Item(935, 9).Value = "Tercera"
$ws.Cells.Item(935, 10).Value = 16000
$ws.Cells.Item(935, 11).Value = 50
$ws.Cells.Item(935, 12).Value = 50
$ws.Cells.Item(935, 13).Value = 50
$ws.Cells.Item(935, 16).Value = 50
$ws.Cells.Item(936, 4).Value = 44340
$ws.Cells.Item(936, 9).Value = "Primera"
$ws.Cells.Item(936, 10).Value = 81500
$ws.Cells.Item(936, 11).Value = 100
$ws.Cells.Item(936, 12).Value = 110
$ws.Cells.Item(936, 13).Value = 106
$ws.Cells.Item(936, 16).Value = 106
$ws.Cells.Item(937, 4).Value = 44340
$ws.Cells.Item(937, 9).Value = "Segunda"
$ws.Cells.Item(937, 10).Value = 43500
$ws.Cells.Item(937, 11).Value = 70
$ws.Cells.Item(937, 12).Value = 85
$ws.Cells.Item(937, 13).Value = 76
$ws.Cells.Item(937, 16).Value = 76
$ws.Cells.Item(938, 4).Value = 44340
$ws.Cells.Item(938, 9).Value = "Tercera"
$ws.Cells.Item(938, 10).Value = 12000
$ws.Cells.Item(938, 11).Value = 50
$ws.Cells.Item(938, 12).Value = 50
$ws.Cells.Item(938, 13).Value = 50
$ws.Cells.Item(938, 16).Value = 50
$ws.Cells.Item(939, 4).Value = 44496
$ws.Cells.Item(939, 9).Value = "Primera"
$ws.Cells.Item(939, 10).Value = 44000
$ws.Cells.Item(939, 11).Value = 90
$ws.Cells.Item(939, 12).Value = 100
$ws.Cells.Item(939, 13).Value = 94
$ws.Cells.Item(939, 16).Value = 94
$ws.Cells.Item(940, 4).Value = 44496
$ws.Cells.Item(940, 9).Value = "Segunda"
$ws.Cells.Item(940, 10).Value = 37000
$ws.Cells.Item(940, 11).Value = 80
$ws.Cells.Item(940, 12).Value = 85
$ws.Cells.Item(940, 13).Value = 82
$ws.Cells.Item(940, 16).Value = 82
$ws.Cells.Item(941, 4).Value = 44496
$ws.Cells.Item(941, 9).Value = "Tercera"
$ws.Cells.Item(941, 10).Value = 9000
$ws.Cells.Item(941, 11).Value = 65
$ws.Cells.Item(941, 12).Value = 65
$ws.Cells.Item(941, 13).Value = 65
$ws.Cells.Item(941, 16).Value = 65
$ws.Cells.Item(942, 4).Value = 44399
$ws.Cells.Item(942, 9).Value = "Primera"
$ws.Cells.Item(942, 10).Value = 56000
$ws.Cells.Item(942, 11).Value = 110
$ws.Cells.Item(942, 12).Value = 125
$ws.Cells.Item(942, 13).Value = 118
$ws.Cells.Item(942, 16).Value = 118
$ws.Cells.Item(943, 4).Value = 44399
$ws.Cells.Item(943, 9).Value = "Segunda"
$ws.Cells.Item(943, 10).Value = 33000
$ws.Cells.Item(943, 11).Value = 80
$ws.Cells.Item(943, 12).Value = 90
$ws.Cells.Item(943, 13).Value = 85
$ws.Cells.Item(943, 16).Value = 85
$ws.Cells.Item(944, 4).Value = 44399
$ws.Cells.Item(944, 9).Value = "Tercera"
$ws.Cells.Item(944, 10).Value = 14000
$ws.Cells.Item(944, 11).Value = 60
$ws.Cells.Item(944, 12).Value = 60
$ws.Cells.Item(944, 13).Value = 60
$ws.Cells.Item(944, 16).Value = 60
$ws.Cells.Item(945, 4).Value = 44425
$ws.Cells.Item(945, 9).Value = "Primera"
$ws.Cells.Item(945, 10).Value = 60000
$ws.Cells.Item(945, 11).Value = 110
$ws.Cells.Item(945, 12).Value = 120
$ws.Cells.Item(945, 13).Value = 114
$ws.Cells.Item(945, 16).Value = 114
$ws.Cells.Item(946, 4).Value = 44425
$ws.Cells.Item(946, 9).Value = "Segunda"
$ws.Cells.Item(946, 10).Value = 51000
$ws.Cells.Item(946, 11).Value = 80
$ws.Cells.Item(946, 12).Value = 95
$ws.Cells.Item(946, 13).Value = 86
$ws.Cells.Item(946, 16).Value = 86
$ws.Cells.Item(947, 4).Value = 44425
$ws.Cells.Item(947, 9).Value = "Tercera"
$ws.Cells.Item(947, 10).Value = 16000
$ws.Cells.Item(947, 11).Value = 70
$ws.Cells.Item(947, 12).Value = 70
$ws.Cells.Item(947, 13).Value = 70
$ws.Cells.Item(947, 16).Value = 70
$ws.Cells.Item(948, 4).Value = 44377
$ws.Cells.Item(948, 9).Value = "Primera"
$ws.Cells.Item(948, 10).Value = 35000
$ws.Cells.Item(948, 11).Value = 85
$ws.Cells.Item(948, 12).Value = 90
$ws.Cells.Item(948, 13).Value = 87
$ws.Cells.Item(948, 16).Value = 87
$ws.Cells.Item(949, 4).Value = 44377
$ws.Cells.Item(949, 9).Value = "Segunda"
$ws.Cells.Item(949, 10).Value = 27000
$ws.Cells.Item(949, 11).Value = 70
$ws.Cells.Item(949, 12).Value = 75
$ws.Cells.Item(949, 13).Value = 72
$ws.Cells.Item(949, 16).Value = 72
$ws.Cells.Item(950, 4).Value = 44377
$ws.Cells.Item(950, 9).Value = "Tercera"
$ws.Cells.Item(950, 10).Value = 9000
$ws.Cells.Item(950, 11).Value = 50
$ws.Cells.Item(950, 12).Value = 50
$ws.Cells.Item(950, 13).Value = 50
$ws.Cells.Item(950, 16).Value = 50
$ws.Cells.Item(951, 4).Value = 44512
$ws.Cells.Item(951, 9).Value = "Primera"
$ws.Cells.Item(951, 10).Value = 54000
$ws.Cells.Item(951, 11).Value = 80
$ws.Cells.Item(951, 12).Value = 90
$ws.Cells.Item(951, 13).Value = 84
$ws.Cells.Item(951, 16).Value = 84
$ws.Cells.Item(952, 4).Value = 44512
$ws.Cells.Item(952, 9).Value = "Segunda"
$ws.Cells.Item(952, 10).Value = 44000
$ws.Cells.Item(952, 11).Value = 70
$ws.Cells.Item(952, 12).Value = 75
$ws.Cells.Item(952, 13).Value = 72
$ws.Cells.Item(952, 16).Value = 72
$ws.Cells.Item(953, 4).Value = 44512
$ws.Cells.Item(953, 9).Value = "Tercera"
$ws.Cells.Item(953, 10).Value = 15000
$ws.Cells.Item(953, 11).Value = 50
$ws.Cells.Item(953, 12).Value = 50
$ws.Cells.Item(953, 13).Value = 50
$ws.Cells.Item(953, 16).Value = 50
$ws.Cells.Item(954, 4).Value = 44397
$ws.Cells.Item(954, 9).Value = "Primera"
$ws.Cells.Item(954, 10).Value = 52000
$ws.Cells.Item(954, 11).Value = 100
$ws.Cells.Item(954, 12).Value = 120
$ws.Cells.Item(954, 13).Value = 109
$ws.Cells.Item(954, 16).Value = 109
$ws.Cells.Item(955, 4).Value = 44397
$ws.Cells.Item(955, 9).Value = "Segunda"
$ws.Cells.Item(955, 10).Value = 45000
$ws.Cells.Item(955, 11).Value = 75
$ws.Cells.Item(955, 12).Value = 90
$ws.Cells.Item(955, 13).Value = 83
$ws.Cells.Item(955, 16).Value = 83
$ws.Cells.Item(956, 4).Value = 44397
$ws.Cells.Item(956, 9).Value = "Tercera"
$ws.Cells.Item(956, 10).Value = 13000
$ws.Cells.Item(956, 11).Value = 60
$ws.Cells.Item(956, 12).Value = 60
$ws.Cells.Item(956, 13).Value = 60
$ws.Cells.Item(956, 16).Value = 60
$ws.Cells.Item(957, 4).Value = 44181
$ws.Cells.Item(957, 9).Value = "Primera"
$ws.Cells.Item(957, 10).Value = 25000
$ws.Cells.Item(957, 11).Value = 80
$ws.Cells.Item(957, 12).Value = 90
$ws.Cells.Item(957, 13).Value = 85
$ws.Cells.Item(957, 16).Value = 85
$ws.Cells.Item(958, 4).Value = 44181
$ws.Cells.Item(958, 9).Value = "Segunda"
$ws.Cells.Item(958, 10).Value = 12000
$ws.Cells.Item(958, 11).Value = 70
$ws.Cells.Item(958, 12).Value = 70
$ws.Cells.Item(958, 13).Value = 70
$ws.Cells.Item(958, 16).Value = 70
$ws.Cells.Item(959, 4).Value = 44497
$ws.Cells.Item(959, 9).Value = "Primera"
$ws.Cells.Item(959, 10).Value = 48000
$ws.Cells.Item(959, 11).Value = 90
$ws.Cells.Item(959, 12).Value = 100
$ws.Cells.Item(959, 13).Value = 95
$ws.Cells.Item(959, 16).Value = 95
$ws.Cells.Item(960, 4).Value = 44497
$ws.Cells.Item(960, 9).Value = "Segunda"
$ws.Cells.Item(960, 10).Value = 39000
$ws.Cells.Item(960, 11).Value = 80
$ws.Cells.Item(960, 12).Value = 85
$ws.Cells.Item(960, 13).Value = 82
$ws.Cells.Item(960, 16).Value = 82
$ws.Cells.Item(961, 4).Value = 44497
$ws.Cells.Item(961, 9).Value = "Tercera"
$ws.Cells.Item(961, 10).Value = 13000
$ws.Cells.Item(961, 11).Value = 60
$ws.Cells.Item(961, 12).Value = 60
$ws.Cells.Item(961, 13).Value = 60
$ws.Cells.Item(961, 16).Value = 60
$ws.Cells.Item(962, 4).Value = 44285
$ws.Cells.Item(962, 9).Value = "Primera"
$ws.Cells.Item(962, 10).Value = 37000
$ws.Cells.Item(962, 11).Value = 110
$ws.Cells.Item(962, 12).Value = 130
$ws.Cells.Item(962, 13).Value = 127
$ws.Cells.Item(962, 16).Value = 127
$ws.Cells.Item(963, 4).Value = 44285
$ws.Cells.Item(963, 9).Value = "Segunda"
$ws.Cells.Item(963, 10).Value = 17500
$ws.Cells.Item(963, 11).Value = 80
$ws.Cells.Item(963, 12).Value = 100
$ws.Cells.Item(963, 13).Value = 97
$ws.Cells.Item(963, 16).Value = 97
$ws.Cells.Item(964, 4).Value = 44362
$ws.Cells.Item(964, 9).Value = "Primera"
$ws.Cells.Item(964, 10).Value = 50000
$ws.Cells.Item(964, 11).Value = 95
$ws.Cells.Item(964, 12).Value = 100
$ws.Cells.Item(964, 13).Value = 97
$ws.Cells.Item(964, 16).Value = 97
$ws.Cells.Item(965, 4).Value = 44362
$ws.Cells.Item(965, 9).Value = "Segunda"
$ws.Cells.Item(965, 10).Value = 41000
$ws.Cells.Item(965, 11).Value = 75
$ws.Cells.Item(965, 12).Value = 80
$ws.Cells.Item(965, 13).Value = 77
$ws.Cells.Item(965, 16).Value = 77
$ws.Cells.Item(966, 4).Value = 44362
$ws.Cells.Item(966, 9).Value = "Tercera"
$ws.Cells.Item(966, 10).Value = 8000
$ws.Cells.Item(966, 11).Value = 50
$ws.Cells.Item(966, 12).Value = 50
$ws.Cells.Item(966, 13).Value = 50
$ws.Cells.Item(966, 16).Value = 50
$ws.Cells.Item(967, 4).Value = 44557
$ws.Cells.Item(967, 9).Value = "Primera"
$ws.Cells.Item(967, 10).Value = 54000
$ws.Cells.Item(967, 11).Value = 70
$ws.Cells.Item(967, 12).Value = 80
$ws.Cells.Item(967, 13).Value = 75
$ws.Cells.Item(967, 16).Value = 75
$ws.Cells.Item(968, 4).Value = 44557
$ws.Cells.Item(968, 9).Value = "Segunda"
$ws.Cells.Item(968, 10).Value = 39000
$ws.Cells.Item(968, 11).Value = 60
$ws.Cells.Item(968, 12).Value = 65
$ws.Cells.Item(968, 13).Value = 62
$ws.Cells.Item(968, 16).Value = 62
$ws.Cells.Item(969, 4).Value = 44557
$ws.Cells.Item(969, 9).Value = "Tercera"
$ws.Cells.Item(969, 10).Value = 13000
$ws.Cells.Item(969, 11).Value = 50
$ws.Cells.Item(969, 12).Value = 50
$ws.Cells.Item(969, 13).Value = 50
$ws.Cells.Item(969, 16).Value = 50
$ws.Cells.Item(970, 4).Value = 44747
$ws.Cells.Item(970, 9).Value = "Primera"
$ws.Cells.Item(970, 10).Value = 36000
$ws.Cells.Item(970, 11).Value = 110
$ws.Cells.Item(970, 12).Value = 120
$ws.Cells.Item(970, 13).Value = 115
$ws.Cells.Item(970, 16).Value = 115
$ws.Cells.Item(971, 4).Value = 44747
$ws.Cells.Item(971, 9).Value = "Segunda"
$ws.Cells.Item(971, 10).Value = 29000
$ws.Cells.Item(971, 11).Value = 85
$ws.Cells.Item(971, 12).Value = 95
$ws.Cells.Item(971, 13).Value = 89
$ws.Cells.Item(971, 16).Value = 89
$ws.Cells.Item(972, 4).Value = 44747
$ws.Cells.Item(972, 9).Value = "Tercera"
$ws.Cells.Item(972, 10).Value = 8000
$ws.Cells.Item(972, 11).Value = 70
$ws.Cells.Item(972, 12).Value = 70
$ws.Cells.Item(972, 13).Value = 70
$ws.Cells.Item(972, 16).Value = 70
$ws.Cells.Item(973, 4).Value = 44357
$ws.Cells.Item(973, 9).Value = "Primera"
$ws.Cells.Item(973, 10).Value = 51000
$ws.Cells.Item(973, 11).Value = 90
$ws.Cells.Item(973, 12).Value = 100
$ws.Cells.Item(973, 13).Value = 95
$ws.Cells.Item(973, 16).Value = 95
$ws.Cells.Item(974, 4).Value = 44357
$ws.Cells.Item(974, 9).Value = "Segunda"
$ws.Cells.Item(974, 10).Value = 41000
$ws.Cells.Item(974, 11).Value = 70
$ws.Cells.Item(974, 12).Value = 75
$ws.Cells.Item(974, 13).Value = 72
$ws.Cells.Item(974, 16).Value = 72
$ws.Cells.Item(975, 4).Value = 44357
$ws.Cells.Item(975, 9).Value = "Tercera"
$ws.Cells.Item(975, 10).Value = 11000
$ws.Cells.Item(975, 11).Value = 50
$ws.Cells.Item(975, 12).Value = 50
$ws.Cells.Item(975, 13).Value = 50
$ws.Cells.Item(975, 16).Value = 50
$ws.Cells.Item(976, 4).Value = 44279
$ws.Cells.Item(976, 9).Value = "Primera"
$ws.Cells.Item(976, 10).Value = 59000
$ws.Cells.Item(976, 11).Value = 110
$ws.Cells.Item(976, 12).Value = 130
$ws.Cells.Item(976, 13).Value = 120
$ws.Cells.Item(976, 16).Value = 120
$ws.Cells.Item(977, 4).Value = 44279
$ws.Cells.Item(977, 9).Value = "Segunda"
$ws.Cells.Item(977, 10).Value = 27000
$ws.Cells.Item(977, 11).Value = 100
$ws.Cells.Item(977, 12).Value = 100
$ws.Cells.Item(977, 13).Value = 100
$ws.Cells.Item(977, 16).Value = 100
$ws.Cells.Item(978, 4).Value = 44551
$ws.Cells.Item(978, 9).Value = "Primera"
$ws.Cells.Item(978, 10).Value = 62000
$ws.Cells.Item(978, 11).Value = 70
$ws.Cells.Item(978, 12).Value = 80
$ws.Cells.Item(978, 13).Value = 74
$ws.Cells.Item(978, 16).Value = 74
$ws.Cells.Item(979, 4).Value = 44551
$ws.Cells.Item(979, 9).Value = "Segunda"
$ws.Cells.Item(979, 10).Value = 49000
$ws.Cells.Item(979, 11).Value = 60
$ws.Cells.Item(979, 12).Value = 65
$ws.Cells.Item(979, 13).Value = 62
$ws.Cells.Item(979, 16).Value = 62
$ws.Cells.Item(980, 4).Value = 44551
$ws.Cells.Item(980, 9).Value = "Tercera"
$ws.Cells.Item(980, 10).Value = 46000
$ws.Cells.Item(980, 11).Value = 50
$ws.Cells.Item(980, 12).Value = 50
$ws.Cells.Item(980, 13).Value = 50
$ws.Cells.Item(980, 16).Value = 50
$ws.Cells.Item(981, 4).Value = 44517
$ws.Cells.Item(981, 9).Value = "Primera"
$ws.Cells.Item(981, 10).Value = 46000
$ws.Cells.Item(981, 11).Value = 75
$ws.Cells.Item(981, 12).Value = 80
$ws.Cells.Item(981, 13).Value = 77
$ws.Cells.Item(981, 16).Value = 77
$ws.Cells.Item(982, 4).Value = 44517
$ws.Cells.Item(982, 9).Value = "Segunda"
$ws.Cells.Item(982, 10).Value = 39000
$ws.Cells.Item(982, 11).Value = 55
$ws.Cells.Item(982, 12).Value = 60
$ws.Cells.Item(982, 13).Value = 57
$ws.Cells.Item(982, 16).Value = 57
$ws.Cells.Item(983, 4).Value = 44517
$ws.Cells.Item(983, 9).Value = "Tercera"
$ws.Cells.Item(983, 10).Value = 12000
$ws.Cells.Item(983, 11).Value = 40
$ws.Cells.Item(983, 12).Value = 40
$ws.Cells.Item(983, 13).Value = 40
$ws.Cells.Item(983, 16).Value = 40
$ws.Cells.Item(984, 4).Value = 44757
$ws.Cells.Item(984, 9).Value = "Primera"
$ws.Cells.Item(984, 10).Value = 32000
$ws.Cells.Item(984, 11).Value = 120
$ws.Cells.Item(984, 12).Value = 130
$ws.Cells.Item(984, 13).Value = 125
$ws.Cells.Item(984, 16).Value = 125
$ws.Cells.Item(985, 4).Value = 44757
$ws.Cells.Item(985, 9).Value = "Segunda"
$ws.Cells.Item(985, 10).Value = 14000
$ws.Cells.Item(985, 11).Value = 90
$ws.Cells.Item(985, 12).Value = 90
$ws.Cells.Item(985, 13).Value = 90
$ws.Cells.Item(985, 16).Value = 90
$ws.Cells.Item(986, 4).Value = 44757
$ws.Cells.Item(986, 9).Value = "Tercera"
$ws.Cells.Item(986, 10).Value = 10000
$ws.Cells.Item(986, 11).Value = 70
$ws.Cells.Item(986, 12).Value = 70
$ws.Cells.Item(986, 13).Value = 70
$ws.Cells.Item(986, 16).Value = 70
$ws.Cells.Item(987, 4).Value = 44547
$ws.Cells.Item(987, 9).Value = "Primera"
$ws.Cells.Item(987, 10).Value = 305000
$ws.Cells.Item(987, 11).Value = 70
$ws.Cells.Item(987, 12).Value = 80
$ws.Cells.Item(987, 13).Value = 79
$ws.Cells.Item(987, 16).Value = 79
$ws.Cells.Item(988, 4).Value = 44547
$ws.Cells.Item(988, 9).Value = "Segunda"
$ws.Cells.Item(988, 10).Value = 51000
$ws.Cells.Item(988, 11).Value = 60
$ws.Cells.Item(988, 12).Value = 65
$ws.Cells.Item(988, 13).Value = 62
$ws.Cells.Item(988, 16).Value = 62
$ws.Cells.Item(989, 4).Value = 44547
$ws.Cells.Item(989, 9).Value = "Tercera"
$ws.Cells.Item(989, 10).Value = 15000
$ws.Cells.Item(989, 11).Value = 50
$ws.Cells.Item(989, 12).Value = 50
$ws.Cells.Item(989, 13).Value = 50
$ws.Cells.Item(989, 16).Value = 50
$ws.Cells.Item(990, 4).Value = 44321
$ws.Cells.Item(990, 9).Value = "Primera"
$ws.Cells.Item(990, 10).Value = 26000
$ws.Cells.Item(990, 11).Value = 100
$ws.Cells.Item(990, 12).Value = 100
$ws.Cells.Item(990, 13).Value = 100
$ws.Cells.Item(990, 16).Value = 100
$ws.Cells.Item(991, 4).Value = 44321
$ws.Cells.Item(991, 9).Value = "Segunda"
$ws.Cells.Item(991, 10).Value = 15000
$ws.Cells.Item(991, 11).Value = 80
$ws.Cells.Item(991, 12).Value = 80
$ws.Cells.Item(991, 13).Value = 80
$ws.Cells.Item(991, 16).Value = 80
$ws.Cells.Item(992, 4).Value = 44438
$ws.Cells.Item(992, 9).Value = "Primera"
$ws.Cells.Item(992, 10).Value = 48000
$ws.Cells.Item(992, 11).Value = 110
$ws.Cells.Item(992, 12).Value = 120
$ws.Cells.Item(992, 13).Value = 115
$ws.Cells.Item(992, 16).Value = 115
$ws.Cells.Item(993, 4).Value = 44438
$ws.Cells.Item(993, 9).Value = "Segunda"
$ws.Cells.Item(993, 10).Value = 31000
$ws.Cells.Item(993, 11).Value = 90
$ws.Cells.Item(993, 12).Value = 95
$ws.Cells.Item(993, 13).Value = 92
$ws.Cells.Item(993, 16).Value = 92
$ws.Cells.Item(994, 4).Value = 44438
$ws.Cells.Item(994, 9).Value = "Tercera"
$ws.Cells.Item(994, 10).Value = 9000
$ws.Cells.Item(994, 11).Value = 65
$ws.Cells.Item(994, 12).Value = 65
$ws.Cells.Item(994, 13).Value = 65
$ws.Cells.Item(994, 16).Value = 65
$ws.Cells.Item(995, 4).Value = 44355
$ws.Cells.Item(995, 9).Value = "Primera"
$ws.Cells.Item(995, 10).Value = 58000
$ws.Cells.Item(995, 11).Value = 100
$ws.Cells.Item(995, 12).Value = 110
$ws.Cells.Item(995, 13).Value = 105
$ws.Cells.Item(995, 16).Value = 105
$ws.Cells.Item(996, 4).Value = 44355
$ws.Cells.Item(996, 9).Value = "Segunda"
$ws.Cells.Item(996, 10).Value = 47000
$ws.Cells.Item(996, 11).Value = 75
$ws.Cells.Item(996, 12).Value = 80
$ws.Cells.Item(996, 13).Value = 77
$ws.Cells.Item(996, 16).Value = 77
$ws.Cells.Item(997, 4).Value = 44355
$ws.Cells.Item(997, 9).Value = "Tercera"
$ws.Cells.Item(997, 10).Value = 15000
$ws.Cells.Item(997, 11).Value = 50
$ws.Cells.Item(997, 12).Value = 50
$ws.Cells.Item(997, 13).Value = 50
$ws.Cells.Item(997, 16).Value = 50
$ws.Cells.Item(998, 4).Value = 44391
$ws.Cells.Item(998, 9).Value = "Primera"
$ws.Cells.Item(998, 10).Value = 33000
$ws.Cells.Item(998, 11).Value = 110
$ws.Cells.Item(998, 12).Value = 120
$ws.Cells.Item(998, 13).Value = 115
$ws.Cells.Item(998, 16).Value = 115
$ws.Cells.Item(999, 4).Value = 44391
$ws.Cells.Item(999, 9).Value = "Segunda"
$ws.Cells.Item(999, 10).Value = 25000
$ws.Cells.Item(999, 11).Value = 85
$ws.Cells.Item(999, 12).Value = 90
$ws.Cells.Item(999, 13).Value = 87
$ws.Cells.Item(999, 16).Value = 87
$ws.Cells.Item(1000, 4).Value = 44391
$ws.Cells.Item(1000, 9).Value = "Tercera"
$ws.Cells.Item(1000, 10).Value = 7000
$ws.Cells.Item(1000, 11).Value = 60
$ws.Cells.Item(1000, 12).Value = 60
$ws.Cells.Item(1000, 13).Value = 60
$ws.Cells.Item(1000, 16).Value = 60
$ws.Cells.Item(1001, 4).Value = 44453
$ws.Cells.Item(1001, 9).Value = "Primera"
$ws.Cells.Item(1001, 10).Value = 49000
$ws.Cells.Item(1001, 11).Value = 110
$ws.Cells.Item(1001, 12).Value = 120
$ws.Cells.Item(1001, 13).Value = 114
$ws.Cells.Item(1001, 16).Value = 114
$ws.Cells.Item(1002, 4).Value = 44453
$ws.Cells.Item(1002, 9).Value = "Segunda"
$ws.Cells.Item(1002, 10).Value = 40000
$ws.Cells.Item(1002, 11).Value = 90
$ws.Cells.Item(1002, 12).Value = 100
$ws.Cells.Item(1002, 13).Value = 94
$ws.Cells.Item(1002, 16).Value = 94
$ws.Cells.Item(1003, 4).Value = 44453
$ws.Cells.Item(1003, 9).Value = "Tercera"
$ws.Cells.Item(1003, 10).Value = 11000
$ws.Cells.Item(1003, 11).Value = 70
$ws.Cells.Item(1003, 12).Value = 70
$ws.Cells.Item(1003, 13).Value = 70
$ws.Cells.Item(1003, 16).Value = 70
$ws.Cells.Item(1004, 4).Value = 44186
$ws.Cells.Item(1004, 9).Value = "Primera"
$ws.Cells.Item(1004, 10).Value = 34000
$ws.Cells.Item(1004, 11).Value = 90
$ws.Cells.Item(1004, 12).Value = 100
$ws.Cells.Item(1004, 13).Value = 95
$ws.Cells.Item(1004, 16).Value = 95
$ws.Cells.Item(1005, 4).Value = 44186
$ws.Cells.Item(1005, 9).Value = "Segunda"
$ws.Cells.Item(1005, 10).Value = 16000
$ws.Cells.Item(1005, 11).Value = 80
$ws.Cells.Item(1005, 12).Value = 80
$ws.Cells.Item(1005, 13).Value = 80
$ws.Cells.Item(1005, 16).Value = 80
$ws.Cells.Item(1006, 4).Value = 44189
$ws.Cells.Item(1006, 9).Value = "Primera"
$ws.Cells.Item(1006, 10).Value = 43000
$ws.Cells.Item(1006, 11).Value = 80
$ws.Cells.Item(1006, 12).Value = 90
$ws.Cells.Item(1006, 13).Value = 85
$ws.Cells.Item(1006, 16).Value = 85
$ws.Cells.Item(1007, 4).Value = 44189
$ws.Cells.Item(1007, 9).Value = "Segunda"
$ws.Cells.Item(1007, 10).Value = 16000
$ws.Cells.Item(1007, 11).Value = 70
$ws.Cells.Item(1007, 12).Value = 70
$ws.Cells.Item(1007, 13).Value = 70
$ws.Cells.Item(1007, 16).Value = 70
$ws.Cells.Item(1008, 4).Value = 44609
$ws.Cells.Item(1008, 9).Value = "Primera"
$ws.Cells.Item(1008, 10).Value = 48000
$ws.Cells.Item(1008, 11).Value = 85
$ws.Cells.Item(1008, 12).Value = 90
$ws.Cells.Item(1008, 13).Value = 87
$ws.Cells.Item(1008, 16).Value = 87
$ws.Cells.Item(1009, 4).Value = 44609
$ws.Cells.Item(1009, 9).Value = "Segunda"
$ws.Cells.Item(1009, 10).Value = 41000
$ws.Cells.Item(1009, 11).Value = 70
$ws.Cells.Item(1009, 12).Value = 75
$ws.Cells.Item(1009, 13).Value = 72
$ws.Cells.Item(1009, 16).Value = 72
$ws.Cells.Item(1010, 4).Value = 44609
$ws.Cells.Item(1010, 9).Value = "Tercera"
$ws.Cells.Item(1010, 10).Value = 12000
$ws.Cells.Item(1010, 11).Value = 60
$ws.Cells.Item(1010, 12).Value = 60
$ws.Cells.Item(1010, 13).Value = 60
$ws.Cells.Item(1010, 16).Value = 60
$ws.Cells.Item(1011, 4).Value = 44489
$ws.Cells.Item(1011, 9).Value = "Primera"
$ws.Cells.Item(1011, 10).Value = 39000
$ws.Cells.Item(1011, 11).Value = 90
$ws.Cells.Item(1011, 12).Value = 100
$ws.Cells.Item(1011, 13).Value = 94
$ws.Cells.Item(1011, 16).Value = 94
$ws.Cells.Item(1012, 4).Value = 44489
$ws.Cells.Item(1012, 9).Value = "Segunda"
$ws.Cells.Item(1012, 10).Value = 34000
$ws.Cells.Item(1012, 11).Value = 75
$ws.Cells.Item(1012, 12).Value = 80
$ws.Cells.Item(1012, 13).Value = 77
$ws.Cells.Item(1012, 16).Value = 77
$ws.Cells.Item(1013, 4).Value = 44489
$ws.Cells.Item(1013, 9).Value = "Tercera"
$ws.Cells.Item(1013, 10).Value = 11000
$ws.Cells.Item(1013, 11).Value = 60
$ws.Cells.Item(1013, 12).Value = 60
$ws.Cells.Item(1013, 13).Value = 60
$ws.Cells.Item(1013, 16).Value = 60
$ws.Cells.Item(1014, 4).Value = 44358
$ws.Cells.Item(1014, 9).Value = "Primera"
$ws.Cells.Item(1014, 10).Value = 46000
$ws.Cells.Item(1014, 11).Value = 90
$ws.Cells.Item(1014, 12).Value = 100
$ws.Cells.Item(1014, 13).Value = 94
$ws.Cells.Item(1014, 16).Value = 94
$ws.Cells.Item(1015, 4).Value = 44358
$ws.Cells.Item(1015, 9).Value = "Segunda"
$ws.Cells.Item(1015, 10).Value = 38000
$ws.Cells.Item(1015, 11).Value = 75
$ws.Cells.Item(1015, 12).Value = 85
$ws.Cells.Item(1015, 13).Value = 79
$ws.Cells.Item(1015, 16).Value = 79
$ws.Cells.Item(1016, 4).Value = 44358
$ws.Cells.Item(1016, 9).Value = "Tercera"
$ws.Cells.Item(1016, 10).Value = 21000
$ws.Cells.Item(1016, 11).Value = 50
$ws.Cells.Item(1016, 12).Value = 60
$ws.Cells.Item(1016, 13).Value = 54
$ws.Cells.Item(1016, 16).Value = 54

# --- Append 3 new rows (1017-1019) carried over from the old tail of the table ---
$ws.Cells.Item(1017, 1).Value = 6
$ws.Cells.Item(1017, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1017, 3).Value = "Metropolitana"
$ws.Cells.Item(1017, 4).Value = 44572
$ws.Cells.Item(1017, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1017, 5).Value = 13
$ws.Cells.Item(1017, 6).Value = 100114014
$ws.Cells.Item(1017, 7).Value = "Betarraga"
$ws.Cells.Item(1017, 8).Value = "Sin especificar"
$ws.Cells.Item(1017, 9).Value = "Primera"
$ws.Cells.Item(1017, 10).Value = 69000
$ws.Cells.Item(1017, 11).Value = 70
$ws.Cells.Item(1017, 12).Value = 75
$ws.Cells.Item(1017, 13).Value = 72
$ws.Cells.Item(1017, 14).Value = "$/unidad"
$ws.Cells.Item(1017, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1017, 16).Value = 72
$ws.Cells.Item(1017, 17).Value = 1
$ws.Cells.Item(1017, 18).Value = "Hortaliza"
$ws.Cells.Item(1018, 1).Value = 6
$ws.Cells.Item(1018, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1018, 3).Value = "Metropolitana"
$ws.Cells.Item(1018, 4).Value = 44572
$ws.Cells.Item(1018, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1018, 5).Value = 13
$ws.Cells.Item(1018, 6).Value = 100114014
$ws.Cells.Item(1018, 7).Value = "Betarraga"
$ws.Cells.Item(1018, 8).Value = "Sin especificar"
$ws.Cells.Item(1018, 9).Value = "Segunda"
$ws.Cells.Item(1018, 10).Value = 54000
$ws.Cells.Item(1018, 11).Value = 60
$ws.Cells.Item(1018, 12).Value = 65
$ws.Cells.Item(1018, 13).Value = 62
$ws.Cells.Item(1018, 14).Value = "$/unidad"
$ws.Cells.Item(1018, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1018, 16).Value = 62
$ws.Cells.Item(1018, 17).Value = 1
$ws.Cells.Item(1018, 18).Value = "Hortaliza"
$ws.Cells.Item(1019, 1).Value = 6
$ws.Cells.Item(1019, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1019, 3).Value = "Metropolitana"
$ws.Cells.Item(1019, 4).Value = 44572
$ws.Cells.Item(1019, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1019, 5).Value = 13
$ws.Cells.Item(1019, 6).Value = 100114014
$ws.Cells.Item(1019, 7).Value = "Betarraga"
$ws.Cells.Item(1019, 8).Value = "Sin especificar"
$ws.Cells.Item(1019, 9).Value = "Tercera"
$ws.Cells.Item(1019, 10).Value = 15000
$ws.Cells.Item(1019, 11).Value = 50
$ws.Cells.Item(1019, 12).Value = 50
$ws.Cells.Item(1019, 13).Value = 50
$ws.Cells.Item(1019, 14).Value = "$/unidad"
$ws.Cells.Item(1019, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1019, 16).Value = 50
$ws.Cells.Item(1019, 17).Value = 1
$ws.Cells.Item(1019, 18).Value = "Hortaliza"
